$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B46 was stored as text "2"; fix it to be a real number 2.
$ws.Range("B46").Value = 2

# Append new row 47 with annotation data from Ying Tang.
$ws.Range("A47").Value = "Ying Tang"

# B47 must stay a text value "3" (not a number) - use the leading
# apostrophe so Excel stores it as text, then reset the style so no
# extra formatting is attached to the cell.
$ws.Range("B47").Value = "'3"
$ws.Range("B47").Style = "Normal"

$ws.Range("C47").Value = "无"
$ws.Range("D47").Value = "SMY"
$ws.Range("E47").Value = "THE"
$ws.Range("F47").Value = "18e2478f-5f8b-460a-bbaf-4b86b95999fd"
$ws.Range("G47").Value = "B1IDRdeCW_annotated.xlsx"
$ws.Range("H47").Value = "This paper presents three observations to understand binary network in Courbariaux, Hubara et al. (2016)."
